$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (G1) to the new
# "Save" header cell (H1), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill the new column with 0 values for the existing data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
